# Applies the 2024-01-21 21:52:57 UTC "Updated cryptos list" GitHub Actions
# refresh: updates Price (D) / Volume(1h) (E) figures for most coins, and
# swaps the Stellar <-> LidoDAOToken rows (37/38) back to rank order.
#
# D-column cells whose new text would otherwise be re-parsed by Excel as a
# number (losing a significant trailing zero, e.g. "32.90" -> 32.9) are
# written via a temporary text NumberFormat, then the cell style is reset
# back to "Normal" so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.718.00'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '2.472.78'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '320.55'
$ws.Range("E5").Value = '  +1.24%  '
$ws.Range("D6").Value = '92.24'
$ws.Range("E6").Value = '  -0.71%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -1.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.90'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.47%  '
$ws.Range("D11").Value = '0.0853'
$ws.Range("E11").Value = '  +1.09%  '
$ws.Range("E12").Value = '  -0.83%  '
$ws.Range("D13").Value = '2.853.65'
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("D15").Value = '15.51'
$ws.Range("E15").Value = '  -1.88%  '
$ws.Range("D16").Value = '2.481.12'
$ws.Range("E16").Value = '  +0.63%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.790'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.35%  '
$ws.Range("D18").Value = '41.647.30'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").Value = '6.43'
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("D20").Value = '0.0₃0940'
$ws.Range("E20").Value = '  -1.90%  '
$ws.Range("D21").Value = '71.54'
$ws.Range("E21").Value = '  +0.44%  '
$ws.Range("E22").Value = '  -3.82%  '
$ws.Range("D23").Value = '239.31'
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("E24").Value = '  +0.79%  '
$ws.Range("D25").Value = '1.94'
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").Value = '24.82'
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("E28").Value = '  -1.27%  '
$ws.Range("D30").Value = '36.35'
$ws.Range("E30").Value = '  +1.43%  '
$ws.Range("D31").Value = '154.96'
$ws.Range("E31").Value = '  -0.72%  '
$ws.Range("E32").Value = '  -1.73%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("E34").Value = '  -0.44%  '
$ws.Range("E35").Value = '  -1.04%  '
$ws.Range("D36").Value = '17.09'
$ws.Range("E36").Value = '  -3.17%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.90'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("B38").Value = 'Stellar'
$ws.Range("C38").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D38").Value = '0.116'
$ws.Range("E38").Value = '  +0.91%  '
$ws.Range("E39").Value = '  +1.68%  '
$ws.Range("E40").Value = '  -0.78%  '
$ws.Range("D41").Value = '3.99'
$ws.Range("E41").Value = '  -1.16%  '
$ws.Range("E42").Value = '  -3.58%  '
$ws.Range("D43").Value = '2.002.96'
$ws.Range("E43").Value = '  +1.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.60'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.59%  '
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("D48").Value = '2.732.22'
$ws.Range("E48").Value = '  +1.07%  '
$ws.Range("D49").Value = '97.28'
$ws.Range("E49").Value = '  +0.20%  '
$ws.Range("D50").Value = '75.76'
$ws.Range("E50").Value = '  +3.52%  '
$ws.Range("D51").Value = '66.95'
$ws.Range("E51").Value = '  -0.36%  '
